# Update countries & provincias Spain
# Applies updated case counts for several countries and refreshes the
# "Datos actualizados" timestamp string on the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 14:22"

# Row 17 - Paises Bajos
$ws.Range("B17").Value = 32655
$ws.Range("C17").Value = 1066
$ws.Range("E17").Value = 28721
$ws.Range("G17").Value = 83
$ws.Range("H17").Value = 3684

# Row 18 - Suiza
$ws.Range("B18").Value = 27740
$ws.Range("C18").Value = 336
$ws.Range("E18").Value = 9272

# Row 22 - Austria
$ws.Range("B22").Value = 14689
$ws.Range("C22").Value = 18
$ws.Range("E22").Value = 3745

# Row 24 - Suecia
$ws.Range("B24").Value = 14385
$ws.Range("C24").Value = 563
$ws.Range("E24").Value = 12295
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 1540

# Row 35 - Dinamarca
$ws.Range("B35").Value = 7384
$ws.Range("C35").Value = 142
$ws.Range("D35").Value = 4141
$ws.Range("E35").Value = 2888
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 355

# Row 60
$ws.Range("E60").Value = 1832
$ws.Range("G60").Value = 5
$ws.Range("H60").Value = 62

# Row 64
$ws.Range("B64").Value = 1871
$ws.Range("C64").Value = 39
$ws.Range("D64").Value = 709
$ws.Range("E64").Value = 1115
$ws.Range("F64").Value = 23
$ws.Range("G64").Value = 8
$ws.Range("H64").Value = 47

# Row 67
$ws.Range("B67").Value = 1676
$ws.Range("C67").Value = 61
$ws.Range("D67").Value = 400
$ws.Range("E67").Value = 1259
